$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - plain numeric values
$ws.Range("A2").Value = "14/07/2023"
$ws.Range("B2").Value = 5000
$ws.Range("C2").Value = 5000
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 5000
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 100

# Row 3 - plain numeric values
$ws.Range("A3").Value = "14/07/2023"
$ws.Range("B3").Value = 6000
$ws.Range("C3").Value = 11000
$ws.Range("D3").Value = 6000.6
$ws.Range("E3").Value = 11000.6
$ws.Range("F3").Value = 0.6
$ws.Range("G3").Value = 100.01

# Row 4 - values formatted as text (two decimals), e.g. "R$" style backup values.
# Leading apostrophe forces these numeric-looking strings to be stored as text
# instead of being auto-converted back to numbers.
$ws.Range("A4").Value = "14/07/2023"
$ws.Range("B4").Value = "'6000.60"
$ws.Range("C4").Value = "'17000.60"
$ws.Range("D4").Value = "'5000.00"
$ws.Range("E4").Value = "'16000.60"
$ws.Range("F4").Value = "'1000.00"
$ws.Range("G4").Value = "'94.12"
